# "Regex fixing..." — fix the UWC-3 and UWQ-1 standard-name regex patterns
# stored in column D of Sheet1, and move the active selection to D6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2: wrap the UWC regex in literal double quotes (content unchanged otherwise)
$ws.Range("D2").Value2 = '"UWC\\D*3"'

# D3: fix the UWQ regex - it had an extra escaped backslash, reduce to a single one
$ws.Range("D3").Value2 = 'UWQ\D*1'

# Move/save the active cell selection to D6 (as last edited by the author)
[void]$ws.Range("D6").Select()
